$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Year:<tab>yyyy<tab><tab><tab>..." -> "Year:<tab>mmm yyyy - mmm yyyy<tab>..."
#   The run holding "yyyy" gets its text replaced, and the three <w:tab/> runs
#   that immediately follow it are removed (one bold tab run + two plain tab
#   runs), leaving just the next (already-present) tab run in place.
# ---------------------------------------------------------------------------

$yyyyRange = $d.Content
$foundYear = $yyyyRange.Find.Execute("yyyy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundYear) {
    # Delete the three single-character tab runs that sit right after "yyyy"
    $tabsAfterYear = $d.Range($yyyyRange.End, $yyyyRange.End + 3)
    $tabsAfterYear.Delete()

    # Now replace the "yyyy" text itself (formatting of the single run is kept)
    $yyyyRange2 = $d.Content
    $yyyyRange2.Find.Execute("yyyy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $yyyyRange2.Text = "mmm yyyy" + [char]0x2013 + "mmm yyyy"
    $yyyyRange2.Text = "mmm yyyy " + [char]0x2013 + " mmm yyyy"
}

# ---------------------------------------------------------------------------
# Change 2: "<tab>Room No :<tab><tab>..." -> "<tab>Room No: nnnnn<tab>..."
#   "Room No" run -> "Room No: "
#   " :" run -> "nnnnn"
#   the single <w:tab/> run right after " :" is removed
# ---------------------------------------------------------------------------

$roomRange = $d.Content
$foundRoom = $roomRange.Find.Execute("Room No", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundRoom) {
    $roomEnd = $roomRange.End

    # Delete the tab run that comes right after "Room No :"
    $tabAfterColon = $d.Range($roomEnd + 2, $roomEnd + 3)
    $tabAfterColon.Delete()

    # "Room No" -> "Room No: "
    $roomRange.Text = "Room No: "

    # " :" -> "nnnnn"
    $colonRange = $d.Range($roomEnd, $roomEnd + 2)
    $colonRange.Text = "nnnnn"
}
